$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 - this shifts the existing row 51
# (and everything below it, through the former row 209) down by one,
# so the former row 209 becomes row 210.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new daily price record.
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 44487
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = 100114013
$ws.Cells.Item(51, 7).Value = "Zanahoria"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 600
$ws.Cells.Item(51, 11).Value = 6500
$ws.Cells.Item(51, 12).Value = 7000
$ws.Cells.Item(51, 13).Value = 6750
$ws.Cells.Item(51, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(51, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 16).Value = 338
$ws.Cells.Item(51, 17).Value = 20
$ws.Cells.Item(51, 18).Value = "Hortaliza"
